$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "38.772.23"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.34%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.103.72"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.16%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "228.69"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.25%  "
$ws.Range("E6").Value = "  +0.13%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "62.41"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.63%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("E9").Value = "  +1.90%  "
$ws.Range("E10").Value = "  -0.18%  "
$ws.Range("E11").Value = "  -1.05%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "15.76"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +6.62%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.414.61"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.07%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "22.12"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.55%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.808"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.40%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.52"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.82%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.105.04"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.70%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "38.795.96"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.55%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "71.96"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.87%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.10"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.17%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0842"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.82%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "228.25"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.54%  "
$ws.Range("E23").Value = "  +0.00%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.35"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.98%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.32"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.34%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "171.50"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.98%  "
$ws.Range("E27").Value = "  +1.42%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.139"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +5.61%  "
$ws.Range("E29").Value = "  +4.49%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "19.37"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.50%  "
$ws.Range("E31").Value = "  +7.76%  "
$ws.Range("E32").Value = "  +0.75%  "
$ws.Range("E33").Value = "  +1.71%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.77"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.56%  "
$ws.Range("E35").Value = "  +7.75%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0617"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.73%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.42"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.34%  "
$ws.Range("E38").Value = "  +1.37%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.999"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.17%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.10"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.84%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "102.90"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.74%  "
$ws.Range("E42").Value = "  +3.38%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.537.02"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.21%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "7.90"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.70%  "
$ws.Range("E45").Value = "  +4.20%  "
$ws.Range("E46").Value = "  -1.06%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0909"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.32%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "4.11"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.53%  "
$ws.Range("E49").Value = "  +1.22%  "
$ws.Range("E50").Value = "  -0.27%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.300.72"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.03%  "
